$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.218.85'
$ws.Range("E2").Value = '  +1.23%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.247.69'
$ws.Range("E3").Value = '  +1.03%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.54'
$ws.Range("E5").Value = '  -2.19%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.98'
$ws.Range("E6").Value = '  -2.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.574'
$ws.Range("E7").Value = '  +0.85%  '

$ws.Range("E8").Value = '  +0.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.528'
$ws.Range("E9").Value = '  -1.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.22'
$ws.Range("E10").Value = '  -2.23%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0817'
$ws.Range("E11").Value = '  -0.62%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.24'
$ws.Range("E12").Value = '  -1.78%  '

$ws.Range("E13").Value = '  +0.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.590.44'
$ws.Range("E14").Value = '  +1.02%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.321.38'
$ws.Range("E15").Value = '  +2.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.836'
$ws.Range("E16").Value = '  -0.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.65'
$ws.Range("E17").Value = '  -3.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '44.076.90'
$ws.Range("E18").Value = '  +1.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0972'
$ws.Range("E19").Value = '  +0.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.42'
$ws.Range("E20").Value = '  +1.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.18'
$ws.Range("E21").Value = '  -6.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.52'
$ws.Range("E22").Value = '  +0.33%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.17'
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("E24").Value = '  -0.76%  '

$ws.Range("E25").Value = '  -1.56%  '

$ws.Range("E26").Value = '  +0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.01'
$ws.Range("E27").Value = '  -0.30%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.22'
$ws.Range("E28").Value = '  -0.37%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '37.64'
$ws.Range("E29").Value = '  +3.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.02'
$ws.Range("E30").Value = '  +1.15%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.19'
$ws.Range("E31").Value = '  +1.53%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '153.04'
$ws.Range("E32").Value = '  -4.32%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0803'
$ws.Range("E33").Value = '  -3.09%  '

$ws.Range("E34").Value = '  +3.30%  '

$ws.Range("E35").Value = '  -3.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.121'
$ws.Range("E36").Value = '  +3.43%  '

$ws.Range("E37").Value = '  -0.86%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.76'
$ws.Range("E38").Value = '  -6.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.48'
$ws.Range("E39").Value = '  -2.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.88'
$ws.Range("E40").Value = '  -2.71%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.66'
$ws.Range("E41").Value = '  -5.70%  '

$ws.Range("E42").Value = '  -2.88%  '

$ws.Range("E43").Value = '  +0.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.746.77'
$ws.Range("E44").Value = '  +2.69%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '83.08'
$ws.Range("E45").Value = '  +1.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.192'
$ws.Range("E46").Value = '  -1.61%  '

$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '100.34'
$ws.Range("E47").Value = '  -1.15%  '

$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.97'
$ws.Range("E48").Value = '  -2.67%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.17'
$ws.Range("E49").Value = '  +1.89%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.01'
$ws.Range("E50").Value = '  -2.76%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.57'
$ws.Range("E51").Value = '  -5.06%  '
